$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row of data (row 8) mirroring the existing rows' structure.
$ws.Range("A8").Value = 42604.890439814815
$ws.Range("A8").NumberFormat = "m/d/yy h:mm"

$ws.Range("B8").Value = "Noun"

$ws.Range("C8").Value = 5839
$ws.Range("D8").Value = 2886
$ws.Range("E8").Value = 570
$ws.Range("F8").Value = 49
$ws.Range("G8").Value = 61
$ws.Range("H8").Value = 44
$ws.Range("I8").Value = 54
$ws.Range("J8").Value = 2
$ws.Range("K8").Value = 8
$ws.Range("L8").Value = 20
$ws.Range("M8").Value = 80
